$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before column D (new D & E); existing D:K shift to F:M
$ws.Columns("D:E").Insert()

# Copy number formatting from column F (the old column D, now shifted) into new D:E
# so the new columns inherit the correct per-row style (date format row 7/38/80, numeric style elsewhere).
$ws.Range("F5:F102").Copy()
$ws.Range("D5:E102").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Populate the two new columns (D,E) with the newest two quarters of data for every row.
# A few historical rows (58, 59, 91) also had their D:M figures restated, so those are
# rewritten in full across D:M rather than only D:E.
$ws.Cells.Item(7, 4).Value = 43465
$ws.Cells.Item(7, 5).Value = 43373
$ws.Cells.Item(8, 4).Value = 2140000
$ws.Cells.Item(8, 5).Value = 1828000
$ws.Cells.Item(9, 4).Value = 702000
$ws.Cells.Item(9, 5).Value = 667000
$ws.Cells.Item(10, 4).Value = 1438000
$ws.Cells.Item(10, 5).Value = 1161000
$ws.Cells.Item(12, 4).Value = "NA"
$ws.Cells.Item(12, 5).Value = "NA"
$ws.Cells.Item(13, 4).Value = 0
$ws.Cells.Item(13, 5).Value = 0
$ws.Cells.Item(14, 4).Value = 0
$ws.Cells.Item(14, 5).Value = 0
$ws.Cells.Item(15, 4).Value = 79000
$ws.Cells.Item(15, 5).Value = 81000
$ws.Cells.Item(17, 4).Value = 2139000
$ws.Cells.Item(17, 5).Value = 1764000
$ws.Cells.Item(18, 4).Value = 1000
$ws.Cells.Item(18, 5).Value = 64000
$ws.Cells.Item(20, 4).Value = 889000
$ws.Cells.Item(20, 5).Value = 621000
$ws.Cells.Item(21, 4).Value = 969000
$ws.Cells.Item(21, 5).Value = 766000
$ws.Cells.Item(22, 4).Value = 0
$ws.Cells.Item(22, 5).Value = 0
$ws.Cells.Item(23, 4).Value = 890000
$ws.Cells.Item(23, 5).Value = 685000
$ws.Cells.Item(24, 4).Value = 32000
$ws.Cells.Item(24, 5).Value = -37000
$ws.Cells.Item(25, 4).Value = 0
$ws.Cells.Item(25, 5).Value = 0
$ws.Cells.Item(26, 4).Value = 858000
$ws.Cells.Item(26, 5).Value = 722000
$ws.Cells.Item(27, 4).Value = 250000
$ws.Cells.Item(27, 5).Value = 163000
$ws.Cells.Item(28, 4).Value = 0
$ws.Cells.Item(28, 5).Value = 0
$ws.Cells.Item(29, 4).Value = 0
$ws.Cells.Item(29, 5).Value = 0
$ws.Cells.Item(30, 4).Value = 0
$ws.Cells.Item(30, 5).Value = 0
$ws.Cells.Item(31, 4).Value = 0
$ws.Cells.Item(31, 5).Value = 0
$ws.Cells.Item(32, 4).Value = -889000
$ws.Cells.Item(32, 5).Value = -621000
$ws.Cells.Item(33, 4).Value = 250000
$ws.Cells.Item(33, 5).Value = 163000
$ws.Cells.Item(34, 4).Value = 0
$ws.Cells.Item(34, 5).Value = 0
$ws.Cells.Item(35, 4).Value = 250000
$ws.Cells.Item(35, 5).Value = 163000
$ws.Cells.Item(38, 4).Value = 43465
$ws.Cells.Item(38, 5).Value = 43373
$ws.Cells.Item(41, 4).Value = 3288000
$ws.Cells.Item(41, 5).Value = 2444000
$ws.Cells.Item(42, 4).Value = 0
$ws.Cells.Item(42, 5).Value = 0
$ws.Cells.Item(43, 4).Value = 1255000
$ws.Cells.Item(43, 5).Value = 694000
$ws.Cells.Item(44, 4).Value = 0
$ws.Cells.Item(44, 5).Value = 0
$ws.Cells.Item(45, 4).Value = 1567000
$ws.Cells.Item(45, 5).Value = 1101000
$ws.Cells.Item(46, 4).Value = 6110000
$ws.Cells.Item(46, 5).Value = 4239000
$ws.Cells.Item(47, 4).Value = 24059000
$ws.Cells.Item(47, 5).Value = 23036000
$ws.Cells.Item(48, 4).Value = 87702000
$ws.Cells.Item(48, 5).Value = 80820000
$ws.Cells.Item(49, 4).Value = 2288000
$ws.Cells.Item(49, 5).Value = 2338000
$ws.Cells.Item(50, 4).Value = 0
$ws.Cells.Item(50, 5).Value = 0
$ws.Cells.Item(51, 4).Value = 0
$ws.Cells.Item(51, 5).Value = 0
$ws.Cells.Item(52, 4).Value = 2361000
$ws.Cells.Item(52, 5).Value = 1150000
$ws.Cells.Item(53, 4).Value = 0
$ws.Cells.Item(53, 5).Value = 0
$ws.Cells.Item(54, 4).Value = 122520000
$ws.Cells.Item(54, 5).Value = 111583000
$ws.Cells.Item(57, 4).Value = 2466000
$ws.Cells.Item(57, 5).Value = 1951000
$ws.Cells.Item(58, 4).Value = 6653000
$ws.Cells.Item(58, 5).Value = 5412000
$ws.Cells.Item(58, 6).Value = 6989000
$ws.Cells.Item(58, 7).Value = 7035000
$ws.Cells.Item(58, 8).Value = 6904000
$ws.Cells.Item(58, 9).Value = 5443000
$ws.Cells.Item(58, 10).Value = 4851000
$ws.Cells.Item(58, 11).Value = 4655000
$ws.Cells.Item(58, 12).Value = 5096000
$ws.Cells.Item(58, 13).Value = 4624000
$ws.Cells.Item(59, 4).Value = 1024000
$ws.Cells.Item(59, 5).Value = 1327000
$ws.Cells.Item(59, 6).Value = 1780000
$ws.Cells.Item(59, 7).Value = 2085000
$ws.Cells.Item(59, 8).Value = 1973000
$ws.Cells.Item(59, 9).Value = 780000
$ws.Cells.Item(59, 10).Value = 643000
$ws.Cells.Item(59, 11).Value = 1635000
$ws.Cells.Item(59, 12).Value = 1309000
$ws.Cells.Item(59, 13).Value = 1101000
$ws.Cells.Item(60, 4).Value = 10143000
$ws.Cells.Item(60, 5).Value = 8690000
$ws.Cells.Item(61, 4).Value = 57942000
$ws.Cells.Item(61, 5).Value = 49912000
$ws.Cells.Item(62, 4).Value = 7695000
$ws.Cells.Item(62, 5).Value = 7223000
$ws.Cells.Item(63, 4).Value = 0
$ws.Cells.Item(63, 5).Value = 0
$ws.Cells.Item(64, 4).Value = 0
$ws.Cells.Item(64, 5).Value = 0
$ws.Cells.Item(65, 4).Value = 0
$ws.Cells.Item(65, 5).Value = 0
$ws.Cells.Item(66, 4).Value = 110163000
$ws.Cells.Item(66, 5).Value = 100367000
$ws.Cells.Item(68, 4).Value = 0
$ws.Cells.Item(68, 5).Value = 0
$ws.Cells.Item(69, 4).Value = 0
$ws.Cells.Item(69, 5).Value = 0
$ws.Cells.Item(70, 4).Value = 4000
$ws.Cells.Item(70, 5).Value = 4000
$ws.Cells.Item(71, 4).Value = 0
$ws.Cells.Item(71, 5).Value = 0
$ws.Cells.Item(72, 4).Value = 0
$ws.Cells.Item(72, 5).Value = 0
$ws.Cells.Item(73, 4).Value = 0
$ws.Cells.Item(73, 5).Value = 0
$ws.Cells.Item(74, 4).Value = 0
$ws.Cells.Item(74, 5).Value = 0
$ws.Cells.Item(75, 4).Value = 0
$ws.Cells.Item(75, 5).Value = 0
$ws.Cells.Item(76, 4).Value = 12353000
$ws.Cells.Item(76, 5).Value = 11212000
$ws.Cells.Item(77, 4).Value = 0
$ws.Cells.Item(77, 5).Value = 0
$ws.Cells.Item(80, 4).Value = 43465
$ws.Cells.Item(80, 5).Value = 43373
$ws.Cells.Item(81, 4).Value = 250000
$ws.Cells.Item(81, 5).Value = 163000
$ws.Cells.Item(83, 4).Value = 79000
$ws.Cells.Item(83, 5).Value = 81000
$ws.Cells.Item(84, 4).Value = 0
$ws.Cells.Item(84, 5).Value = 0
$ws.Cells.Item(85, 4).Value = 0
$ws.Cells.Item(85, 5).Value = 0
$ws.Cells.Item(86, 4).Value = 0
$ws.Cells.Item(86, 5).Value = 0
$ws.Cells.Item(87, 4).Value = 0
$ws.Cells.Item(87, 5).Value = 0
$ws.Cells.Item(88, 4).Value = 0
$ws.Cells.Item(88, 5).Value = 0
$ws.Cells.Item(89, 4).Value = 209000
$ws.Cells.Item(89, 5).Value = 369000
$ws.Cells.Item(91, 4).Value = 6289000
$ws.Cells.Item(91, 5).Value = -2389000
$ws.Cells.Item(91, 6).Value = -1577000
$ws.Cells.Item(91, 7).Value = -2815000
$ws.Cells.Item(91, 8).Value = -84000
$ws.Cells.Item(91, 9).Value = -56000
$ws.Cells.Item(91, 10).Value = -72000
$ws.Cells.Item(91, 11).Value = -1799000
$ws.Cells.Item(91, 12).Value = -3796000
$ws.Cells.Item(91, 13).Value = -1415000
$ws.Cells.Item(92, 4).Value = 0
$ws.Cells.Item(92, 5).Value = 0
$ws.Cells.Item(93, 4).Value = 0
$ws.Cells.Item(93, 5).Value = 0
$ws.Cells.Item(94, 4).Value = -4422000
$ws.Cells.Item(94, 5).Value = -666000
$ws.Cells.Item(96, 4).Value = 0
$ws.Cells.Item(96, 5).Value = 0
$ws.Cells.Item(97, 4).Value = 0
$ws.Cells.Item(97, 5).Value = 0
$ws.Cells.Item(98, 4).Value = 0
$ws.Cells.Item(98, 5).Value = 0
$ws.Cells.Item(99, 4).Value = 0
$ws.Cells.Item(99, 5).Value = 0
$ws.Cells.Item(100, 4).Value = 5043000
$ws.Cells.Item(100, 5).Value = 1157000
$ws.Cells.Item(101, 4).Value = 14000
$ws.Cells.Item(101, 5).Value = -16000
$ws.Cells.Item(102, 4).Value = 844000
$ws.Cells.Item(102, 5).Value = 844000
